$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Cells.Item(2, 2).Value = 0.7596339113680154
$ws.Cells.Item(3, 2).Value = 0.847297768612078
$ws.Cells.Item(4, 2).Value = 0.7596339113680154
$ws.Cells.Item(5, 2).Value = 0.7549191201777731

# --- Sheet: Class Metrics ---
$ws = $wb.Worksheets.Item("Class Metrics")
$ws.Cells.Item(2, 2).Value = 0.9223300970873787
$ws.Cells.Item(2, 3).Value = 0.9405940594059405
$ws.Cells.Item(2, 4).Value = 0.9313725490196079
$ws.Cells.Item(3, 2).Value = 0.8690476190476191
$ws.Cells.Item(3, 3).Value = 0.9798657718120806
$ws.Cells.Item(3, 4).Value = 0.9211356466876972
$ws.Cells.Item(4, 2).Value = 0.7518796992481203
$ws.Cells.Item(4, 4).Value = 0.8583690987124464
$ws.Cells.Item(5, 2).Value = 0.9047619047619048
$ws.Cells.Item(5, 3).Value = 0.76
$ws.Cells.Item(5, 4).Value = 0.8260869565217391
$ws.Cells.Item(6, 2).Value = 0.8
$ws.Cells.Item(6, 3).Value = 0.5
$ws.Cells.Item(6, 4).Value = 0.6153846153846154
$ws.Cells.Item(7, 2).Value = 0.9904761904761905
$ws.Cells.Item(7, 3).Value = 0.485981308411215
$ws.Cells.Item(7, 4).Value = 0.6520376175548589
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = 0.06
$ws.Cells.Item(8, 4).Value = 0.1132075471698113
$ws.Cells.Item(9, 2).Value = 0.8019323671497585
$ws.Cells.Item(9, 3).Value = 0.8645833333333334
$ws.Cells.Item(9, 4).Value = 0.8320802005012531
$ws.Cells.Item(10, 2).Value = 0.6041666666666666
$ws.Cells.Item(10, 3).Value = 0.90625
$ws.Cells.Item(10, 4).Value = 0.725
$ws.Cells.Item(11, 2).Value = 0.5342019543973942
$ws.Cells.Item(11, 3).Value = 0.9213483146067416
$ws.Cells.Item(11, 4).Value = 0.6762886597938145
$ws.Cells.Item(12, 2).Value = 0.8675496688741722
$ws.Cells.Item(12, 3).Value = 0.7797619047619048
$ws.Cells.Item(12, 4).Value = 0.8213166144200627
$ws.Cells.Item(13, 2).Value = 0.8429752066115702
$ws.Cells.Item(13, 3).Value = 0.7183098591549296
$ws.Cells.Item(13, 4).Value = 0.7756653992395437
$ws.Cells.Item(14, 2).Value = 0.9949748743718593
$ws.Cells.Item(14, 3).Value = 0.6149068322981367
$ws.Cells.Item(14, 4).Value = 0.7600767754318618
$ws.Cells.Item(15, 2).Value = 0.2303030303030303
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 0.374384236453202
$ws.Cells.Item(16, 2).Value = 0.9017341040462428
$ws.Cells.Item(16, 3).Value = 0.975
$ws.Cells.Item(16, 4).Value = 0.9369369369369369

# --- Sheet: Confusion Matrix ---
$ws = $wb.Worksheets.Item("Confusion Matrix")
$ws.Cells.Item(2, 2).Value = 95
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 3
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(3, 3).Value = 146
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 76
$ws.Cells.Item(5, 9).Value = 13
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(6, 3).Value = 4
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 8
$ws.Cells.Item(6, 16).Value = 2
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(7, 3).Value = 4
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 7).Value = 104
$ws.Cells.Item(7, 9).Value = 3
$ws.Cells.Item(7, 10).Value = 27
$ws.Cells.Item(7, 11).Value = 58
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 12
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 15).Value = 2
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 4).Value = 13
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 8).Value = 6
$ws.Cells.Item(8, 9).Value = 17
$ws.Cells.Item(8, 10).Value = 16
$ws.Cells.Item(8, 11).Value = 44
$ws.Cells.Item(8, 13).Value = 0
$ws.Cells.Item(8, 15).Value = 2
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 4
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 166
$ws.Cells.Item(9, 10).Value = 9
$ws.Cells.Item(9, 11).Value = 5
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(9, 14).Value = 1
$ws.Cells.Item(9, 16).Value = 2
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = 87
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1
$ws.Cells.Item(10, 15).Value = 4
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 6
$ws.Cells.Item(11, 7).Value = 1
$ws.Cells.Item(11, 9).Value = 4
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 164
$ws.Cells.Item(11, 15).Value = 2
$ws.Cells.Item(12, 3).Value = 7
$ws.Cells.Item(12, 12).Value = 131
$ws.Cells.Item(12, 13).Value = 5
$ws.Cells.Item(12, 15).Value = 19
$ws.Cells.Item(12, 16).Value = 5
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 10).Value = 1
$ws.Cells.Item(13, 11).Value = 19
$ws.Cells.Item(13, 12).Value = 8
$ws.Cells.Item(13, 13).Value = 102
$ws.Cells.Item(13, 16).Value = 8
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 10).Value = 2
$ws.Cells.Item(14, 11).Value = 13
$ws.Cells.Item(14, 12).Value = 10
$ws.Cells.Item(14, 13).Value = 1
$ws.Cells.Item(14, 14).Value = 198
$ws.Cells.Item(14, 15).Value = 97
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 15).Value = 38
$ws.Cells.Item(16, 9).Value = 3
$ws.Cells.Item(16, 15).Value = 1
$ws.Cells.Item(16, 16).Value = 156
